$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Exact "16×78=" "91×18="
Replace-Exact "23×93=" "43×12="
Replace-Exact "50×38=" "57×91="
Replace-Exact "14×26=" "48×14="
Replace-Exact "40×62=" "79×74="
Replace-Exact "68×99=" "12×64="
Replace-Exact "98×26=" "19×94="
Replace-Exact "86×62=" "93×93="
Replace-Exact "56×63=" "51×13="
Replace-Exact "40×76=" "30×42="
Replace-Exact "69×42=" "58×95="
Replace-Exact "11×65=" "58×21="
Replace-Exact "60×54=" "75×61="
Replace-Exact "95×69=" "69×94="
Replace-Exact "15×97=" "94×78="
Replace-Exact "39×17=" "31×23="
Replace-Exact "42×77=" "85×45="
Replace-Exact "18×98=" "36×21="
Replace-Exact "26×95=" "34×35="
Replace-Exact "90×72=" "33×95="
Replace-Exact "52×12=" "20×45="
Replace-Exact "63×53=" "75×33="
Replace-Exact "35×62=" "82×82="
Replace-Exact "92×39=" "16×18="
Replace-Exact "64×22=" "27×27="
